# Simulated Wild Card round and logged it
# Appends per-play/per-kick results to the running season logs and
# updates the aggregate totals on each summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: per-play yardage logs (R = rush, P = pass) for OFF (B) / DEF (C)
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + " 2 2 0 4 22 -5 2 3 3 4 2 3 18 8 -1 8 10 46 1 22 27 3 6 0 4 2 6 1 21 0 11 4 4 3 1"
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + " 3 19 13 11 10 6 0 28 7 1 5 1 6 11 8 12"
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + " 1 1 5 9 -1 2 2 -6 1 1 2 0 7 14 1 -1 3 17 -1 0 2 4 3 0 2 0 15"
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + " 4 11 10 7 0 7 25 2 5 16 8 2 10 10 7 8 -2 8 20 11 3 6 10 6 3 11 14 1 11"

# ---------------------------------------------------------------------------
# OFF sheet: Home (row 2) / Road (row 3) offensive totals
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value2 = 496
$wsOFF.Range("F2").Value2 = 143
$wsOFF.Range("G2").Value2 = 152
$wsOFF.Range("I2").Value2 = 27
$wsOFF.Range("J2").Value2 = 121
$wsOFF.Range("L2").Value2 = 530
$wsOFF.Range("M2").Value2 = 339
$wsOFF.Range("O2").Value2 = 48
$wsOFF.Range("P2").Value2 = 30
$wsOFF.Range("Q2").Value2 = 1144

$wsOFF.Range("B3").Value2 = 20
$wsOFF.Range("C3").Value2 = 352
$wsOFF.Range("D3").Value2 = 9
$wsOFF.Range("E3").Value2 = 65
$wsOFF.Range("F3").Value2 = 240
$wsOFF.Range("H3").Value2 = 60
$wsOFF.Range("I3").Value2 = 124
$wsOFF.Range("J3").Value2 = 71
$wsOFF.Range("N3").Value2 = 40

# ---------------------------------------------------------------------------
# DEF sheet: Home (row 2) / Road (row 3) defensive totals
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value2 = 405
$wsDEF.Range("F2").Value2 = 104
$wsDEF.Range("G2").Value2 = 101
$wsDEF.Range("J2").Value2 = 42
$wsDEF.Range("L2").Value2 = 607
$wsDEF.Range("M2").Value2 = 398
$wsDEF.Range("O2").Value2 = 54
$wsDEF.Range("P2").Value2 = 28
$wsDEF.Range("Q2").Value2 = 1038

$wsDEF.Range("C3").Value2 = 396
$wsDEF.Range("E3").Value2 = 78
$wsDEF.Range("F3").Value2 = 247
$wsDEF.Range("G3").Value2 = 72
$wsDEF.Range("H3").Value2 = 66
$wsDEF.Range("I3").Value2 = 132
$wsDEF.Range("J3").Value2 = 111
$wsDEF.Range("N3").Value2 = 38

# ---------------------------------------------------------------------------
# ST sheet: special-teams totals (row 2) + per-kick logs (rows 3-6)
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value2 = 189
$wsST.Range("D2").Value2 = 130
$wsST.Range("F2").Value2 = 386
$wsST.Range("G2").Value2 = 382
$wsST.Range("J2").Value2 = 181
$wsST.Range("K2").Value2 = 178
$wsST.Range("L2").Value2 = 111
$wsST.Range("M2").Value2 = 100

$wsST.Range("B3").Value2 = 116

$wsST.Range("B4").Value2 = $wsST.Range("B4").Value2 + " 59 52 58"
$wsST.Range("B5").Value2 = $wsST.Range("B5").Value2 + " 32 13 21"
$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + " 19 23 27 26 18"
$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + " 44 40 36 23 50 39 36"
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + " 7 0 0 0 23 0 0"
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + " 0 0 0 0 0 0 0 0"

# ---------------------------------------------------------------------------
# TURNS sheet: Home (row 2) / Road (row 3) turnover totals
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B2").Value2 = 19
$wsTURNS.Range("C2").Value2 = 10
$wsTURNS.Range("D2").Value2 = 20

$wsTURNS.Range("E3").Value2 = 18

# ---------------------------------------------------------------------------
# PEN sheet: OFF penalty counts (column B)
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value2 = 41
$wsPEN.Range("B3").Value2 = 47
